$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Igfbp4", "Fzd8", "ECs", 3, 1, 11.24958966666667, 33.748769, 0.04815412300202451, 0.04815412300202451, 2, 0.6666666666666666, 2.219727, 6.659181, 0.1284523376175785, 0.1284523376175785, 24.971017922021, 224.739161298189, 0.006185509665534454, 0.006185509665534454),
    @("ECs", "Igfbp4", "Fzd8", "FAPs", 3, 1, 11.24958966666667, 33.748769, 0.04815412300202451, 0.04815412300202451, 3, 1, 11.50585933333333, 34.517578, 0.6658271614778302, 0.6658271614778303, 129.4361962623869, 1164.925766361482, 0.03206232303189227, 0.03206232303189227),
    @("ECs", "Igfbp4", "Fzd8", "M2", 3, 1, 11.24958966666667, 33.748769, 0.04815412300202451, 0.04815412300202451, 1, 0.3333333333333333, 0.01066666666666667, 0.032, 0.0006172643158013742, 0.0006172643158013743, 0.1199956231111111, 1.079960608, 0.00002972382178785988, 0.00002972382178785988),
    @("ECs", "Igfbp4", "Fzd8", "sCs", 3, 1, 11.24958966666667, 33.748769, 0.04815412300202451, 0.04815412300202451, 3, 1, 3.544296666666666, 10.63289, 0.2051032365887898, 0.2051032365887898, 39.87188315693444, 358.84694841241, 0.009876566482809919, 0.009876566482809919),
    @("FAPs", "Igfbp4", "Fzd8", "ECs", 3, 1, 208.2711843333333, 624.813553, 0.8915095150431698, 0.8915095150431697, 2, 0.6666666666666666, 2.219727, 6.659181, 0.1284523376175785, 0.1284523376175785, 462.305171186677, 4160.746540680093, 0.1145164812156089, 0.1145164812156089),
    @("FAPs", "Igfbp4", "Fzd8", "FAPs", 3, 1, 208.2711843333333, 624.813553, 0.8915095150431698, 0.8915095150431697, 3, 1, 11.50585933333333, 34.517578, 0.6658271614778302, 0.6658271614778303, 2396.33895012607, 21567.05055113463, 0.5935912498316707, 0.5935912498316707),
    @("FAPs", "Igfbp4", "Fzd8", "M2", 3, 1, 208.2711843333333, 624.813553, 0.8915095150431698, 0.8915095150431697, 1, 0.3333333333333333, 0.01066666666666667, 0.032, 0.0006172643158013742, 0.0006172643158013743, 2.221559299555555, 19.994033696, 0.0005502970108335371, 0.0005502970108335371),
    @("FAPs", "Igfbp4", "Fzd8", "sCs", 3, 1, 208.2711843333333, 624.813553, 0.8915095150431698, 0.8915095150431697, 3, 1, 3.544296666666666, 10.63289, 0.2051032365887898, 0.2051032365887898, 738.1748643953521, 6643.573779558169, 0.1828514869850565, 0.1828514869850565),
    @("M2", "Igfbp4", "Fzd8", "ECs", 3, 1, 3.262842333333333, 9.788527, 0.01396667040408609, 0.01396667040408609, 2, 0.6666666666666666, 2.219727, 6.659181, 0.1284523376175785, 0.1284523376175785, 7.242619224043001, 65.183573016387, 0.001794051462139108, 0.001794051462139107),
    @("M2", "Igfbp4", "Fzd8", "FAPs", 3, 1, 3.262842333333333, 9.788527, 0.01396667040408609, 0.01396667040408609, 3, 1, 11.50585933333333, 34.517578, 0.6658271614778302, 0.6658271614778303, 37.54180491417844, 337.876244227606, 0.009299388510449058, 0.00929938851044906),
    @("M2", "Igfbp4", "Fzd8", "M2", 3, 1, 3.262842333333333, 9.788527, 0.01396667040408609, 0.01396667040408609, 1, 0.3333333333333333, 0.01066666666666667, 0.032, 0.0006172643158013742, 0.0006172643158013743, 0.03480365155555556, 0.313232864, 0.0000086211272510015, 0.000008621127251001502),
    @("M2", "Igfbp4", "Fzd8", "sCs", 3, 1, 3.262842333333333, 9.788527, 0.01396667040408609, 0.01396667040408609, 3, 1, 3.544296666666666, 10.63289, 0.2051032365887898, 0.2051032365887898, 11.56448120589222, 104.08033085303, 0.002864609304246917, 0.002864609304246917),
    @("sCs", "Igfbp4", "Fzd8", "ECs", 3, 1, 10.83271733333333, 32.498152, 0.04636969155071965, 0.04636969155071963, 2, 0.6666666666666666, 2.219727, 6.659181, 0.1284523376175785, 0.1284523376175785, 24.045675148168, 216.411076333512, 0.005956295274296017, 0.005956295274296015),
    @("sCs", "Igfbp4", "Fzd8", "FAPs", 3, 1, 10.83271733333333, 32.498152, 0.04636969155071965, 0.04636969155071963, 3, 1, 11.50585933333333, 34.517578, 0.6658271614778302, 0.6658271614778303, 124.6397218350951, 1121.757496515856, 0.03087420010381819, 0.03087420010381818),
    @("sCs", "Igfbp4", "Fzd8", "M2", 3, 1, 10.83271733333333, 32.498152, 0.04636969155071965, 0.04636969155071963, 1, 0.3333333333333333, 0.01066666666666667, 0.032, 0.0006172643158013742, 0.0006172643158013743, 0.1155489848888889, 1.039940864, 0.00002862235592897572, 0.00002862235592897572),
    @("sCs", "Igfbp4", "Fzd8", "sCs", 3, 1, 10.83271733333333, 32.498152, 0.04636969155071965, 0.04636969155071963, 3, 1, 3.544296666666666, 10.63289, 0.2051032365887898, 0.2051032365887898, 38.39436393547555, 345.54927541928, 0.009510573816676459, 0.009510573816676457)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowData[$j]
    }
}

Write-Host "Rows written:" $rows.Length
